$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values are assigned in the same order the original author typed them
# so that newly-created shared-string entries land at the same indices as
# the authoritative workbook.

# Row 11
$ws.Cells.Item(11,3).Value = "PO3_DGW_SIQ_08"
$ws.Cells.Item(11,5).Value = "Change mode, don't save adjustments"
$ws.Cells.Item(11,7).Value = "21/2/2020"
$ws.Cells.Item(11,4).Value = "If MODE button is pressed during time adjust or alarm adjust, does it:`n1. Change mode without saving.`n2. Change mode with saving.`n3. Stay in the same mode without saving (Cancel). `n4. Stay in the same mode without saving?"

# Row 12 (question/answer)
$ws.Cells.Item(12,4).Value = "What is displayed on LCD in Alarm mode?"
$ws.Cells.Item(12,5).Value = "Current time, until alarm time is reached "

# Row 13 question
$ws.Cells.Item(13,4).Value = "If Alarm time is adjusted, and mode is changed to display time or stopwatch, does it still trigger buzzer if alarm time is reached?"

# SIQ IDs for rows 12 / 13
$ws.Cells.Item(12,3).Value = "PO3_DGW_SIQ_09"
$ws.Cells.Item(13,3).Value = "PO3_DGW_SIQ_10"

# Row 14 question
$ws.Cells.Item(14,4).Value = "If alarm is working in all modes, How can we set alarm to infinity? (Never trigger a buzzer)"

# SIQ IDs for rows 14 / 15
$ws.Cells.Item(14,3).Value = "PO3_DGW_SIQ_11"
$ws.Cells.Item(15,3).Value = "PO3_DGW_SIQ_12"

# Row 15 question
$ws.Cells.Item(15,4).Value = "What are the initial states for current time, alarm time?"

# Row 14 / 15 answers
$ws.Cells.Item(14,5).Value = "Set to a time that's not valid (e.g: 00:00:00)"
$ws.Cells.Item(15,5).Value = "Current time: 12:00:00 AM`nAlarm time: infinity (e.g: 00:00:00)"

# Row 16 question + SIQ id
$ws.Cells.Item(16,4).Value = "What are the limits for incrementing in adjust mode (for both: current time and alarm)? What happens if we increment to a boundary value?"
$ws.Cells.Item(16,3).Value = "PO3_DGW_SIQ_13"

# Requirement IDs
$ws.Cells.Item(12,1).Value = "PO3_DGW_CRS_F_02"
$ws.Cells.Item(15,1).Value = "PO3_DGW_CRS_F_01`nPO3_DGW_CRS_F_02"
$ws.Cells.Item(16,1).Value = "PO3_DGW_CRS_KE_03"

# Row 16 rich-text answer (bold labels "Current time:" / "Alarm: ")
$c16 = $ws.Cells.Item(16,5)
$c16.Value = "Current time:" + [char]10 + "Hours -> 1:12 (reset to 1)" + [char]10 + "Minutes -> 0:59 (reset to 0)" + [char]10 + "keep AM/PM as they were if hours are incremented above boundary value." + [char]10 + "Alarm: " + "Same, but add an option for infinity"
$c16.Characters(1,13).Font.Bold = $true
$c16.Characters(1,13).Font.Size = 12
$c16.Characters(1,13).Font.Name = "Calibri"
$c16.Characters(14,128).Font.Size = 12
$c16.Characters(14,128).Font.Name = "Calibri"
$c16.Characters(142,7).Font.Bold = $true
$c16.Characters(142,7).Font.Size = 12
$c16.Characters(142,7).Font.Name = "Calibri"
$c16.Characters(149,36).Font.Size = 12
$c16.Characters(149,36).Font.Name = "Calibri"

# Remaining cells that reuse already-existing shared strings
$ws.Cells.Item(11,1).Value = "General"
$ws.Cells.Item(13,1).Value = "PO3_DGW_CRS_F_02"
$ws.Cells.Item(14,1).Value = "PO3_DGW_CRS_F_02"

$ws.Cells.Item(11,2).Value = "Mariam"
$ws.Cells.Item(12,2).Value = "Mariam"
$ws.Cells.Item(13,2).Value = "Mariam"
$ws.Cells.Item(14,2).Value = "Mariam"
$ws.Cells.Item(15,2).Value = "Mariam"
$ws.Cells.Item(16,2).Value = "Mariam"

$ws.Cells.Item(11,6).Value = "Mariam"
$ws.Cells.Item(12,6).Value = "Mariam"
$ws.Cells.Item(14,6).Value = "Mariam"
$ws.Cells.Item(15,6).Value = "Mariam"
$ws.Cells.Item(16,6).Value = "Mariam"

$ws.Cells.Item(12,7).Value = "21/2/2020"
$ws.Cells.Item(13,7).Value = "21/2/2020"
$ws.Cells.Item(14,7).Value = "21/2/2020"
$ws.Cells.Item(15,7).Value = "21/2/2020"
$ws.Cells.Item(16,7).Value = "21/2/2020"

$ws.Cells.Item(11,8).Value = "Not answered"
$ws.Cells.Item(12,8).Value = "Not answered"
$ws.Cells.Item(13,8).Value = "Not answered"
$ws.Cells.Item(14,8).Value = "Not answered"
$ws.Cells.Item(15,8).Value = "Not answered"
$ws.Cells.Item(16,8).Value = "Not answered"

# Row heights
$ws.Rows.Item(11).RowHeight = 120
$ws.Rows.Item(12).RowHeight = 79.95
$ws.Rows.Item(13).RowHeight = 79.95
$ws.Rows.Item(14).RowHeight = 79.95
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 150

# Selection, matching the saved view state in the edited workbook
$ws.Range("D12").Select()
